# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Fri Mar 24 19:50:10 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.825.04'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '1.764.42'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4260'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3636'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.44%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07518'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.33%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.095'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9996'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.77'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.077'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.290'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('D16').Value = '1.780.17'
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001062'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06383'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.914'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.37%  '
$ws.Range('D23').Value = '27.864.49'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.120'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.50'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').Value = '1.989.69'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.154'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.40'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.120'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.681'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.564'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08885'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02294'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2108'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06039'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6342'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.975'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.178'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9992'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.897'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.394'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5880'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.690'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.987'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.189'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06835'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.11%  '
